$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 32: a new "Counting" passive skill entry
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = "Counting"
$ws.Range("C32").Value = "PassiveSkill"
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0

# Match style/format of the row above (row 31)
$ws.Range("A31:E31").Copy()
$ws.Range("A32:E32").PasteSpecial(-4122)  # xlPasteFormats

# Re-apply values since PasteSpecial formats only (values preserved anyway, but keep explicit)
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = "Counting"
$ws.Range("C32").Value = "PassiveSkill"
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0

# Update selection/view to match new last row
$ws.Range("A32").Select()
